$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "45+47=92"
$t.Cell(1,2).Range.Text = "85-58=27"
$t.Cell(1,3).Range.Text = "5+79=84"
$t.Cell(1,4).Range.Text = "71-23=48"
$t.Cell(1,5).Range.Text = "85+8=93"
$t.Cell(2,1).Range.Text = "66-7=59"
$t.Cell(2,2).Range.Text = "55-17=38"
$t.Cell(2,3).Range.Text = "24-15=9"
$t.Cell(2,4).Range.Text = "59+22=81"
$t.Cell(2,5).Range.Text = "87-68=19"
$t.Cell(3,1).Range.Text = "61-35=26"
$t.Cell(3,2).Range.Text = "73-25=48"
$t.Cell(3,3).Range.Text = "5+57=62"
$t.Cell(3,4).Range.Text = "67-48=19"
$t.Cell(3,5).Range.Text = "18+46=64"
$t.Cell(4,1).Range.Text = "71-6=65"
$t.Cell(4,2).Range.Text = "90-13=77"
$t.Cell(4,3).Range.Text = "13+29=42"
$t.Cell(4,4).Range.Text = "18+76=94"
$t.Cell(4,5).Range.Text = "23+38=61"
$t.Cell(5,1).Range.Text = "67+17=84"
$t.Cell(5,2).Range.Text = "39+43=82"
$t.Cell(5,3).Range.Text = "57+34=91"
$t.Cell(5,4).Range.Text = "27+29=56"
$t.Cell(5,5).Range.Text = "68+17=85"
$t.Cell(6,1).Range.Text = "36+25=61"
$t.Cell(6,2).Range.Text = "9+29=38"
$t.Cell(6,3).Range.Text = "9+69=78"
$t.Cell(6,4).Range.Text = "82-6=76"
$t.Cell(6,5).Range.Text = "41-17=24"
$t.Cell(7,1).Range.Text = "80-2=78"
$t.Cell(7,2).Range.Text = "12+19=31"
$t.Cell(7,3).Range.Text = "32-14=18"
$t.Cell(7,4).Range.Text = "58+4=62"
$t.Cell(7,5).Range.Text = "26+18=44"
$t.Cell(8,1).Range.Text = "7+58=65"
$t.Cell(8,2).Range.Text = "40-7=33"
$t.Cell(8,3).Range.Text = "54+29=83"
$t.Cell(8,4).Range.Text = "86-48=38"
$t.Cell(8,5).Range.Text = "23-9=14"
$t.Cell(9,1).Range.Text = "71-35=36"
$t.Cell(9,2).Range.Text = "83-27=56"
$t.Cell(9,3).Range.Text = "51-28=23"
$t.Cell(9,4).Range.Text = "68+28=96"
$t.Cell(9,5).Range.Text = "50-12=38"
$t.Cell(10,1).Range.Text = "29+24=53"
$t.Cell(10,2).Range.Text = "48+24=72"
$t.Cell(10,3).Range.Text = "9+73=82"
$t.Cell(10,4).Range.Text = "5+69=74"
$t.Cell(10,5).Range.Text = "25+19=44"
$t.Cell(11,1).Range.Text = "19+42=61"
$t.Cell(11,2).Range.Text = "80-25=55"
$t.Cell(11,3).Range.Text = "62-13=49"
$t.Cell(11,4).Range.Text = "19+2=21"
$t.Cell(11,5).Range.Text = "56+6=62"
$t.Cell(12,1).Range.Text = "52-16=36"
$t.Cell(12,2).Range.Text = "45+29=74"
$t.Cell(12,3).Range.Text = "7+5=12"
$t.Cell(12,4).Range.Text = "54+17=71"
$t.Cell(12,5).Range.Text = "47+47=94"
$t.Cell(13,1).Range.Text = "91-87=4"
$t.Cell(13,2).Range.Text = "27+6=33"
$t.Cell(13,3).Range.Text = "22+19=41"
$t.Cell(13,4).Range.Text = "82-75=7"
$t.Cell(13,5).Range.Text = "51-19=32"
$t.Cell(14,1).Range.Text = "71-38=33"
$t.Cell(14,2).Range.Text = "65-57=8"
$t.Cell(14,3).Range.Text = "94-8=86"
$t.Cell(14,4).Range.Text = "95-19=76"
$t.Cell(14,5).Range.Text = "58+5=63"
$t.Cell(15,1).Range.Text = "96-49=47"
$t.Cell(15,2).Range.Text = "28+33=61"
$t.Cell(15,3).Range.Text = "37+58=95"
$t.Cell(15,4).Range.Text = "18+17=35"
$t.Cell(15,5).Range.Text = "77+19=96"
$t.Cell(16,1).Range.Text = "24+48=72"
$t.Cell(16,2).Range.Text = "70-2=68"
$t.Cell(16,3).Range.Text = "18+33=51"
$t.Cell(16,4).Range.Text = "91-9=82"
$t.Cell(16,5).Range.Text = "72-14=58"
$t.Cell(17,1).Range.Text = "31-16=15"
$t.Cell(17,2).Range.Text = "84-36=48"
$t.Cell(17,3).Range.Text = "29+32=61"
$t.Cell(17,4).Range.Text = "62-54=8"
$t.Cell(17,5).Range.Text = "40-8=32"
$t.Cell(18,1).Range.Text = "26+65=91"
$t.Cell(18,2).Range.Text = "83-35=48"
$t.Cell(18,3).Range.Text = "38+55=93"
$t.Cell(18,4).Range.Text = "42+29=71"
$t.Cell(18,5).Range.Text = "2+79=81"
$t.Cell(19,1).Range.Text = "36-17=19"
$t.Cell(19,2).Range.Text = "60-22=38"
$t.Cell(19,3).Range.Text = "13+18=31"
$t.Cell(19,4).Range.Text = "41-33=8"
$t.Cell(19,5).Range.Text = "25+16=41"
$t.Cell(20,1).Range.Text = "33+38=71"
$t.Cell(20,2).Range.Text = "70-15=55"
$t.Cell(20,3).Range.Text = "72-14=58"
$t.Cell(20,4).Range.Text = "17+6=23"
$t.Cell(20,5).Range.Text = "62+29=91"
